# Update workbook to reflect data through 2022-08-10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name / sheet name)
$ws.Name = "Through 2022-08-10"

# Update the column header label for the 2022 column (cell I1)
$ws.Range("I1").Value = "2022 (through 08-10)"

# Update the August 2022 value (row 9 = August)
$ws.Range("I9").Value = 52

# Update the Total 2022 value (row 14 = Total)
$ws.Range("I14").Value = 1022
